$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 24
$ws.Range("F3").Value = 480
$ws.Range("F5").Value = 81
$ws.Range("F6").Value = 324
$ws.Range("F7").Value = 1308
$ws.Range("F8").Value = 516
$ws.Range("F9").Value = 107
$ws.Range("F10").Value = 1327
$ws.Range("F13").Value = 137
$ws.Range("F16").Value = 117
$ws.Range("F17").Value = 253
$ws.Range("F18").Value = 1679
$ws.Range("F19").Value = 624
$ws.Range("F21").Value = 248
$ws.Range("F22").Value = 2738
$ws.Range("F23").Value = 24
$ws.Range("F24").Value = 406
$ws.Range("F26").Value = 931
$ws.Range("F27").Value = 1217
$ws.Range("F29").Value = 2842
$ws.Range("F30").Value = 1643
$ws.Range("F31").Value = 86
$ws.Range("F32").Value = 122
$ws.Range("F33").Value = 684
$ws.Range("F34").Value = 871
$ws.Range("F35").Value = 1883
$ws.Range("F37").Value = 1887
$ws.Range("F39").Value = 25
$ws.Range("F40").Value = 14
$ws.Range("F42").Value = 47
$ws.Range("F44").Value = 808
$ws.Range("F45").Value = 1038
$ws.Range("F46").Value = 115
$ws.Range("F47").Value = 445
$ws.Range("F48").Value = 226
$ws.Range("F49").Value = 3364

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 188
$ws.Range("F12").Value = 807
$ws.Range("F13").Value = 25
$ws.Range("F17").Value = 11

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 480
$ws.Range("F4").Value = 81
$ws.Range("F6").Value = 188
$ws.Range("F7").Value = 324
$ws.Range("F8").Value = 1308
$ws.Range("F9").Value = 516
$ws.Range("F10").Value = 107
$ws.Range("F11").Value = 1327
$ws.Range("F14").Value = 137
$ws.Range("F17").Value = 117
$ws.Range("F18").Value = 253
$ws.Range("F19").Value = 1679
$ws.Range("F20").Value = 624
$ws.Range("F22").Value = 248
$ws.Range("F23").Value = 2739
$ws.Range("F24").Value = 24
$ws.Range("F25").Value = 406
$ws.Range("F27").Value = 1217
$ws.Range("F28").Value = 2842
$ws.Range("F29").Value = 1643
$ws.Range("F30").Value = 86
$ws.Range("F32").Value = 122
$ws.Range("F33").Value = 807
$ws.Range("F34").Value = 25
$ws.Range("F35").Value = 871
$ws.Range("F36").Value = 1883
$ws.Range("F37").Value = 11
$ws.Range("F39").Value = 1887
$ws.Range("F42").Value = 808
$ws.Range("F43").Value = 1038
$ws.Range("F44").Value = 115
$ws.Range("F45").Value = 445
$ws.Range("F47").Value = 226
$ws.Range("F48").Value = 3364
